$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report gained one new line item ("ROSASEEF 30 CAPS") inserted as the
# new row 32 ("م" = 26), pushing every row from the old row 32 down through
# the totals/footer rows down by one. We replicate that by copying full rows
# (values + formats) from bottom to top so we don't disturb existing styles,
# then fill the freed row 32 with the new item, and refresh the computed
# total and the generated timestamp string.

for ($r = 51; $r -ge 32; $r--) {
    $destRow = $r + 1
    $src = $ws.Range("A" + $r + ":Q" + $r)
    $dst = $ws.Range("A" + $destRow + ":Q" + $destRow)
    $src.Copy()
    $dst.PasteSpecial(-4104) # xlPasteAll
}
$excel.CutCopyMode = 0

# Fill in the new row 32 with the new product line
$ws.Range("A32").Value2 = 26
$ws.Range("C32").Value2 = "ROSASEEF 30 CAPS"
$ws.Range("H32").Value2 = "0:0"
$ws.Range("L32").Value2 = 1
$ws.Range("N32").Value2 = "120.00"
$ws.Range("P32").Value2 = "39.6000"
$ws.Range("Q32").Value2 = "0:1"

# Renumber the "م" (index) column for the shifted rows so it stays 26..44
$n = 27
for ($r = 33; $r -le 50; $r++) {
    $ws.Range("A" + $r).Value2 = $n
    $n = $n + 1
}

# Update the running total (row 51 now, was row 50) to include the new line
$ws.Range("P51").Value2 = 1888.925

# Refresh the generation timestamp shown in the footer (row 52 now, was 51)
$ws.Range("A52").Value2 = "Sunday, 7 September, 2025 5:47 PM"
